# Update the "想去人数" (interested-count) values on the "展览" and
# "全部类型" worksheets, matching the values produced by the latest
# gh-pages data generation run.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 95
    $ws.Range("F7").Value = 133
    $ws.Range("F9").Value = 357
}
